# The slide master currently carries the "Integral" (Red Violet) theme
# colors. Re-apply the default Office theme's 12-slot color scheme to the
# slide master's theme, matching the built-in "Office Theme" palette.
#
# MsoThemeColorSchemeIndex order for ThemeColorScheme.Colors():
#   1 dk1   2 lt1   3 dk2   4 lt2
#   5 accent1  6 accent2  7 accent3  8 accent4  9 accent5  10 accent6
#   11 hlink   12 folHlink
#
# Note: .RGB uses the standard OLE BGR-packed integer (0x00BBGGRR), so each
# hex RRGGBB value below is byte-reversed before assignment.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$tcs = $master.Theme.ThemeColorScheme

$tcs.Colors(1).RGB  = 0x000000   # dk1      000000
$tcs.Colors(2).RGB  = 0xFFFFFF   # lt1      FFFFFF
$tcs.Colors(3).RGB  = 0x6A5444   # dk2      44546A
$tcs.Colors(4).RGB  = 0xE6E6E7   # lt2      E7E6E6
$tcs.Colors(5).RGB  = 0xD59B5B   # accent1  5B9BD5
$tcs.Colors(6).RGB  = 0x317DED   # accent2  ED7D31
$tcs.Colors(7).RGB  = 0xA5A5A5   # accent3  A5A5A5
$tcs.Colors(8).RGB  = 0x00C0FF   # accent4  FFC000
$tcs.Colors(9).RGB  = 0xC47244   # accent5  4472C4
$tcs.Colors(10).RGB = 0x47AD70   # accent6  70AD47
$tcs.Colors(11).RGB = 0xC16305   # hlink    0563C1
$tcs.Colors(12).RGB = 0x724F95   # folHlink 954F72
